$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: seed the shared-string table with the new unique strings in the
# exact order they need to appear (so their shared-string indices come out
# as 32..48, matching the target file). We stage them in a scratch column
# far away from the used range, then clear the scratch cells again.
$newStrings = @(
    "Department of Justice",
    "Gayunman,",
    "Hiniling",
    "Inaasahang,",
    "Isinagawa ng DOJ",
    "Isinagawa ng NBI",
    "Ito",
    "Kasabay",
    "Nagbabala",
    "National Bureau of Investigation",
    "Nauna",
    "PAOCTF",
    "Pangulo",
    "Presidential Anti-Organized Crime Task Force",
    "Sa",
    "Tukoy",
    "Ulat"
)

$scratchCol = 50
for ($i = 0; $i -lt $newStrings.Length; $i++) {
    $ws.Cells.Item($i + 1, $scratchCol).Value = $newStrings[$i]
}
for ($i = 0; $i -lt $newStrings.Length; $i++) {
    $ws.Cells.Item($i + 1, $scratchCol).ClearContents()
}

# --- Step 2: rewrite column A (PERSON entries) and column B (the constant
# "PERSON" label) for rows 5-35 with the final, alphabetically-resorted
# list that now also includes the newly-recognised entities.
$finalA = @(
    "Dacer at Corbito",
    "Dacer",
    "Department of Justice",
    "Ebdane",
    "Emmanuel Corbito",
    "Gayunman,",
    "Grace Amargo at Joy Cantos",
    "Hiniling",
    "Inaasahang,",
    "Isinagawa ng DOJ",
    "Isinagawa ng NBI",
    "Ito",
    "Jimmy Lopez at Alex Diloy,",
    "Kasabay",
    "Lopez at Diloy",
    "NBI Director Reynaldo Wycoco",
    "NBI-National Capital Region Director Samuel Ong",
    "Nagbabala",
    "National Bureau of Investigation",
    "Nauna",
    "Ong",
    "PAOCTF Chief Director Hermogenes Ebdane",
    "PAOCTF",
    "PAOCTF-Visayas Chief Sr. Supt. Teofilo Vina",
    "Pangulo",
    "Pangulong Gloria Macapagal-Arroyo",
    "Presidential Anti-Organized Crime Task Force",
    "Sa",
    "Salvador Bubby"" Dacer""",
    "Tukoy",
    "Ulat"
)

$startRow = 5
for ($i = 0; $i -lt $finalA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $finalA[$i]
    $ws.Cells.Item($r, 2).Value = "PERSON"
}

# --- Step 3: update the sheet selection to match the new used range.
$ws.Range("A1:B35").Select()
